$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18-22 down to 19-23
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with the new weekly price entry
$ws.Range("A18").Value = 4
$ws.Range("B18").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C18").Value = "Los Lagos"
$ws.Range("D18").Value = 44782
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 100112013
$ws.Range("G18").Value = "Alcachofa"
$ws.Range("H18").Value = "Madrigal"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 90
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 15000
$ws.Range("M18").Value = 15000
$ws.Range("N18").Value = "$/caja 40 unidades"
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 375
$ws.Range("Q18").Value = 40
$ws.Range("R18").Value = "Hortaliza"
